$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain numeric-looking strings as TEXT
# (e.g. "120.50", "1.00", "0.130"). Force text format on the cells we
# touch so Excel does not auto-coerce them into numbers and strip
# significant trailing zeros / reformat them.
$touchedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $touchedRows) {
    $ws.Range("D$r`:E$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "43.634.47"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "2.271.88"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "120.50"
$ws.Range("E5").Value = "  +7.92%  "

$ws.Range("D6").Value = "264.98"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +3.39%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +2.61%  "

$ws.Range("D10").Value = "47.53"
$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").Value = "9.41"
$ws.Range("E12").Value = "  +5.80%  "

$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("D14").Value = "15.60"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").Value = "0.895"
$ws.Range("E15").Value = "  +4.82%  "

$ws.Range("D16").Value = "2.613.87"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "2.270.75"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "43.598.85"
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").Value = "0.0000109"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("D20").Value = "6.93"
$ws.Range("E20").Value = "  +0.95%  "

$ws.Range("D21").Value = "72.14"
$ws.Range("E21").Value = "  +1.30%  "

$ws.Range("D22").Value = "2.40"
$ws.Range("E22").Value = "  -3.23%  "

$ws.Range("D23").Value = "234.65"
$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "2.88"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "9.44"
$ws.Range("E25").Value = "  -2.36%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "12.08"
$ws.Range("E26").Value = "  +6.18%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.02"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").Value = "41.84"
$ws.Range("E28").Value = "  +3.46%  "

$ws.Range("D29").Value = "3.35"
$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "173.97"
$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("D32").Value = "21.41"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").Value = "0.0918"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  -1.23%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.0384"
$ws.Range("E36").Value = "  +9.27%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "4.19"
$ws.Range("E37").Value = "  +8.63%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "4.59"
$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.110"
$ws.Range("E39").Value = "  +5.17%  "

$ws.Range("D40").Value = "2.54"
$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("D41").Value = "13.78"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").Value = "72.67"
$ws.Range("E42").Value = "  -3.88%  "

$ws.Range("D43").Value = "0.238"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("D46").Value = "5.79"
$ws.Range("E46").Value = "  -5.14%  "

$ws.Range("D47").Value = "77.24"
$ws.Range("E47").Value = "  +46.97%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "102.67"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.655"
$ws.Range("E50").Value = "  +15.45%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "8.49"
$ws.Range("E51").Value = "  -2.38%  "
